# Update the cryptos list with refreshed prices / 1h volume percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (D value, E value). A key is omitted when that
# column did not change for the row (row 5 has no changes at all).
$updates = @{
    2  = @{ D = "19.960.87"; E = "  -5.06%  " }
    3  = @{ D = "1.415.87";  E = "  -5.70%  " }
    4  = @{            E = "  -0.63%  " }
    6  = @{ D = "276.29";    E = "  -2.51%  " }
    7  = @{ D = "0.3665";    E = "  -3.90%  " }
    8  = @{ D = "0.3095";    E = "  -0.61%  " }
    9  = @{ D = "39.71";     E = "  -7.40%  " }
    10 = @{ D = "1.034";     E = "  +0.16%  " }
    11 = @{ D = "0.06534";   E = "  -5.34%  " }
    12 = @{            E = "  -0.70%  " }
    13 = @{ D = "5.499";     E = "  -1.14%  " }
    14 = @{ D = "17.63";     E = "  -0.37%  " }
    15 = @{ D = "6.201";     E = "  -2.38%  " }
    16 = @{ D = "1.416.33";  E = "  -5.97%  " }
    17 = @{ D = "0.00001020"; E = "  -4.06%  " }
    18 = @{ D = "0.05667";   E = "  -13.46%  " }
    19 = @{            E = "  -0.58%  " }
    20 = @{ D = "71.25";     E = "  -12.88%  " }
    21 = @{ D = "5.621";     E = "  -5.49%  " }
    22 = @{ D = "14.74";     E = "  -1.90%  " }
    23 = @{ D = "10.92";     E = "  +0.79%  " }
    24 = @{ D = "2.248";     E = "  -3.84%  " }
    25 = @{ D = "19.978.14"; E = "  -4.99%  " }
    26 = @{ D = "2.270";     E = "  -2.04%  " }
    27 = @{ D = "133.04";    E = "  -9.36%  " }
    28 = @{ D = "17.29";     E = "  -3.42%  " }
    29 = @{ D = "1.576.98";  E = "  -5.94%  " }
    30 = @{ D = "109.89";    E = "  -3.62%  " }
    31 = @{ D = "3.887";     E = "  -19.00%  " }
    32 = @{ D = "5.270";     E = "  -9.40%  " }
    33 = @{ D = "0.8195";    E = "  -13.47%  " }
    34 = @{ D = "0.07701";   E = "  -2.43%  " }
    35 = @{            E = "  +0.96%  " }
    36 = @{ D = "8.337";     E = "  -0.74%  " }
    37 = @{ D = "4.925";     E = "  -2.28%  " }
    38 = @{ D = "0.05793";   E = "  +0.67%  " }
    39 = @{ D = "1.000";     E = "  -0.55%  " }
    40 = @{ D = "0.02053";   E = "  -2.98%  " }
    41 = @{ D = "10.49";     E = "  -5.33%  " }
    42 = @{ D = "0.1888";    E = "  -3.81%  " }
    43 = @{ D = "1.100";     E = "  -5.09%  " }
    44 = @{ D = "12.41";     E = "  -3.14%  " }
    45 = @{ D = "0.5316";    E = "  -4.77%  " }
    46 = @{            E = "  -3.24%  " }
    47 = @{ D = "0.5192";    E = "  -3.93%  " }
    48 = @{ D = "115.91";    E = "  +2.69%  " }
    49 = @{ D = "1.772";     E = "  -3.38%  " }
    50 = @{ D = "1.034";     E = "  -8.25%  " }
    51 = @{ D = "1.000";     E = "  -0.63%  " }
}

# Rows whose new "Price" text contains two '.' separators (e.g. "1.415.87")
# already fail Excel's numeric parse, so they stay text without any extra
# help. Everything else looks like a plain decimal ("276.29") and needs to
# be pinned to text explicitly or Excel will silently coerce it to a number.
$twoDotRows = @(2, 3, 16, 25, 29)

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        if ($twoDotRows -notcontains $row) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}
